# Update "想去人数" (interested count) figures for the
# "合肥·第九届环形宇宙动漫游戏嘉年华" and "合肥·MAX特摄同人only2.0" events.
# These rows are duplicated on the "展览" sheet (rows 3-4) and the
# "全部类型" sheet (rows 7-8).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2710
$ws1.Range("F4").Value = 122

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2710
$ws4.Range("F8").Value = 122
